# Adds a new week of price data (Primera / Segunda) for Betarraga at
# "Terminal La Palmera de La Serena" as of date serial 45013 (2023-03-28).
# The new rows are inserted at the top of the existing weekly data block
# (current rows 400-401), pushing the rest of the historical rows down by
# two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 400, shifting all the
# existing data (old rows 400:427) down to 402:429.
$ws.Rows.Item(400).Resize(2).Insert()

# --- New row 400: Betarraga, calidad "Primera" ---
$ws.Range("A400").Value = 8
$ws.Range("B400").Value = "Terminal La Palmera de La Serena"
$ws.Range("C400").Value = "Coquimbo"
$ws.Range("D400").Value = 45013
$ws.Range("E400").Value = 4
$ws.Range("F400").Value = 100114014
$ws.Range("G400").Value = "Betarraga"
$ws.Range("H400").Value = "Sin especificar"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 1900
$ws.Range("K400").Value = 500
$ws.Range("L400").Value = 600
$ws.Range("M400").Value = 550
$ws.Range("N400").Value = "`$/paquete 3 unidades"
$ws.Range("O400").Value = "Provincia del Elquí"
$ws.Range("P400").Value = 183
$ws.Range("Q400").Value = 3
$ws.Range("R400").Value = "Hortaliza"

# --- New row 401: Betarraga, calidad "Segunda" ---
$ws.Range("A401").Value = 8
$ws.Range("B401").Value = "Terminal La Palmera de La Serena"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 45013
$ws.Range("E401").Value = 4
$ws.Range("F401").Value = 100114014
$ws.Range("G401").Value = "Betarraga"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Segunda"
$ws.Range("J401").Value = 1260
$ws.Range("K401").Value = 400
$ws.Range("L401").Value = 450
$ws.Range("M401").Value = 425
$ws.Range("N401").Value = "`$/paquete 3 unidades"
$ws.Range("O401").Value = "Provincia del Elquí"
$ws.Range("P401").Value = 142
$ws.Range("Q401").Value = 3
$ws.Range("R401").Value = "Hortaliza"
